$wb = $excel.ActiveWorkbook

# Rename sheets: the sheet that "doesn't work" is fixed, becomes "cases".
# The old "cases" sheet becomes the fallback "caseswithoutmaize".
$wsBroken = $wb.Worksheets.Item("cases_doesn't work")
$wsOldCases = $wb.Worksheets.Item("cases")

$wsOldCases.Name = "caseswithoutmaize"
$wsBroken.Name = "cases"

# Remove the now-obsolete warning note row (row 7) from the fixed "cases" sheet.
$wsBroken.Rows.Item(7).Delete()

# Make the fixed "cases" sheet the active/selected sheet & cell.
$wsBroken.Select()
$wsBroken.Range("B5").Select()

# Restore previous "cases" sheet selection (no longer active tab).
$wsOldCases.Range("F9").Select()
$wsBroken.Activate()
